# Applies the "handles float input without breaking stuff" marksheet update:
#  - summary block (rows 10-12) gets real grading numbers instead of the
#    "not graded yet" placeholders
#  - C11 becomes a genuine numeric -1 (was stored as text "-1")
#  - per-question student answers (columns A and, for the first three
#    questions, D) are now filled in and colour-coded correct/incorrect
#    against the adjacent "Correct Ans" column
#  - the unused G/H "3rd answer pair" block and the now-redundant D/E block
#    (rows 19-40) are cleared out, which also shrinks the sheet's used range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Copy-Style($fromAddr, $toAddr) {
    $ws.Range($fromAddr).Copy()
    $ws.Range($toAddr).PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = $false
}

# --- Row 10-12 labels pick up the "mtitleStyle" look (same style as the
#     A9 header row) now that the table is fully populated ---
Copy-Style "A9" "A10"
Copy-Style "A9" "A11"
Copy-Style "A9" "A12"

# --- Row 10: Right / Wrong / Not Attempted / Max ---
$ws.Range("B10").Value = 24
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 28

# --- Row 11: per-question marking scheme (now a real number, not text) ---
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# --- Row 12: totals ---
$ws.Range("B12").Value = 96
$ws.Range("C12").Value = -1
$ws.Range("E12").Value = "95/112"

# --- Drop the unused third Student/Correct-Ans block (columns G:H) ---
$ws.Range("G15:H40").Clear()

# --- Per-question answers in column A (student) vs column B (correct).
#     Rows 20, 21 and 25 were left unanswered, so column A stays blank. ---
$answersA = @{
    16 = "Option B"
    17 = "Option D"
    18 = "Option B"
    19 = "Option C"
    22 = "Option D"
    23 = "Option D"
    24 = "Option A"
    26 = "Option C"
    27 = "Option A"
    28 = "Option D"
    29 = "Option D"
    30 = "Option B"
    31 = "Option D"
    32 = "Option C"
    33 = "Option D"
    34 = "Option B"
    35 = "Option D"
    36 = "Option A"
    37 = "Option A"
    38 = "Option A"
    39 = "Option D"
    40 = "Option D"
}

foreach ($row in $answersA.Keys) {
    $cellAddr = "A$row"
    $studentAns = $answersA[$row]
    $correctAns = $ws.Range("B$row").Value2

    $ws.Range($cellAddr).Value = $studentAns
    if ($studentAns -eq $correctAns) {
        Copy-Style "B10" $cellAddr   # B10 carries the "correctStyle" (green)
    } else {
        Copy-Style "C10" $cellAddr   # C10 carries the "incorrectStyle" (red)
    }
}

# --- The first three questions also had a second Student/Correct-Ans pair
#     in columns D:E; fill those in (all answered correctly) ---
$answersD = @{
    16 = "Option A"
    17 = "Option C"
    18 = "Option D"
}
foreach ($row in $answersD.Keys) {
    $cellAddr = "D$row"
    $studentAns = $answersD[$row]
    $correctAns = $ws.Range("E$row").Value2

    $ws.Range($cellAddr).Value = $studentAns
    if ($studentAns -eq $correctAns) {
        Copy-Style "B10" $cellAddr
    } else {
        Copy-Style "C10" $cellAddr
    }
}

# --- The D:E pair is only needed for rows 16-18; everything below that was
#     a duplicate of the per-question grid and is no longer used ---
$ws.Range("D19:E40").Clear()

Write-Output "done"
